# The models were retrained.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update accuracy/loss values for rows 2-4 with retrained model results
$ws.Range("E2").Value = 0.6542
$ws.Range("E3").Value = 0.4798
$ws.Range("F3").Value = 0.8221
$ws.Range("E4").Value = 0.3349
$ws.Range("F4").Value = 0.8981

# Clear out the stale results for rows 5-11 (pending retraining / not yet available)
$ws.Range("E5:F11").ClearContents()

# Update the active selection to reflect where the user left off
$ws.Range("F4").Select()
